# Regenerate merged AHB files:
#  - rename the "_old" / "_new" header-row column suffixes to the actual
#    merged form-version tags ("_FV2210" / "_FV2304")
#  - wrap the data range in a table ("Table1")
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

# Columns 1-10 (A-J): "<Header>_old"  -> "<Header>_FV2210"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i] + "_FV2210"
}

# Column 11 (K): "diff" stays as-is.

# Columns 12-21 (L-U): "<Header>_new" -> "<Header>_FV2304"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $headers[$i] + "_FV2304"
}

# Turn the used range into a proper table.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U67"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
